# Applies the changes described by the commit diff to the Trello work-report
# workbook. All edits are plain cell-value updates (numbers, text labels,
# and a few brand-new cells / one removed cell) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section "Aktiv" (rows 17-21)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 13
$ws.Range("G17").Value = 11

# C18 / D19 hold digit-strings stored as TEXT in the original file. Prefix
# with an apostrophe so Excel keeps them as text instead of coercing them to
# numbers.
$ws.Range("C18").Value = "'10"
$ws.Range("D19").Value = "'18"

$ws.Range("G19").Value = 2
$ws.Range("D20").Value = "38.5% der Karten"
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = 11

# ---------------------------------------------------------------------------
# Section "Karten ohne Aktivität" (rows 27-31)
# ---------------------------------------------------------------------------
$ws.Range("B27").Value = "Revisión de Mapas y Material para Sistema de  Riego por Goteo y Material"
$ws.Range("C27").Value = 11

$ws.Range("B28").Value = "NEU: Gewächshäuser"
$ws.Range("C28").Value = 10

$ws.Range("B29").Value = "Log Frame erstellen"
$ws.Range("C29").Value = 9

$ws.Range("C30").Value = 9

$ws.Range("B31").Value = "Workshop - Nutrición Escolar ⭐️"
$ws.Range("C31").Value = 6

# ---------------------------------------------------------------------------
# Section "Gemeinschaftlich" - Aktivste Mitglieder (rows 34-43)
# ---------------------------------------------------------------------------
$ws.Range("F34").Value = 5
$ws.Range("G34").Value = "(38.5%)"

$ws.Range("B35").Value = "Leonie Ziller"
$ws.Range("C35").Value = 6

$ws.Range("B36").Value = "Maria Lütticke"
$ws.Range("C36").Value = 4

$ws.Range("B37").Value = "Jonas Ullmann"
$ws.Range("C37").Value = 2

$ws.Range("B38").Value = "Chantal Bußmann"
$ws.Range("C38").Value = 2

$ws.Range("B39").Value = "Micha Landoll"
$ws.Range("C39").Value = 2

$ws.Range("E39").Value = "Leonie Ziller"
$ws.Range("F39").Value = 6

$ws.Range("E40").Value = "Maria Lütticke"
$ws.Range("F40").Value = 4

$ws.Range("E41").Value = "Jonas Ullmann"
$ws.Range("F41").Value = 2

$ws.Range("E42").Value = "Chantal Bußmann"
$ws.Range("F42").Value = 2

$ws.Range("E43").Value = "Eva Greven"
# F43 stays at 1 (unchanged)

# ---------------------------------------------------------------------------
# Section "Zuverlässig" - Meiste Abgeschlossene Karten / Aktivste Helfer
# (rows 48-52) - these rows gain brand-new name cells alongside the count
# updates.
# ---------------------------------------------------------------------------
$ws.Range("B48").Value = "Jonas Ullmann"
$ws.Range("C48").Value = 1

$ws.Range("B49").Value = "Leonie Ziller"
$ws.Range("C49").Value = 1

$ws.Range("B50").Value = "Maria Lütticke"
$ws.Range("C50").Value = 1
$ws.Range("E50").Value = "Wafic Sabbagh"
$ws.Range("F50").Value = 1

$ws.Range("B51").Value = "Chantal Bußmann"
$ws.Range("C51").Value = 1
$ws.Range("E51").Value = "Micha Landoll"
$ws.Range("F51").Value = 1

$ws.Range("B52").Value = "Micha Landoll"
$ws.Range("C52").Value = 1
# F52 stays at 0 (unchanged)

# ---------------------------------------------------------------------------
# Section "Pünktlich" (rows 60-64)
# ---------------------------------------------------------------------------
$ws.Range("A60").Value = "Leonie Ziller"
$ws.Range("B60").Value = 2
$ws.Range("G60").Value = 2

$ws.Range("A61").Value = "Jonas Ullmann"
# B61 stays at 1 (unchanged)
$ws.Range("G61").Value = 3

$ws.Range("A62").Value = "Maria Lütticke"
$ws.Range("B62").Value = 1

$ws.Range("A63").Value = "Chantal Bußmann"
$ws.Range("B63").Value = 1

$ws.Range("A64").Value = "Micha Landoll"
$ws.Range("B64").Value = 1

# ---------------------------------------------------------------------------
# Section "Detailliert" (rows 70-74)
# ---------------------------------------------------------------------------
$ws.Range("B70").Value = 13
$ws.Range("F70").Value = "Leonie Ziller"
$ws.Range("G70").Value = 4

$ws.Range("B71").Value = 4
$ws.Range("F71").Value = "Maria Lütticke"
$ws.Range("G71").Value = 3

$ws.Range("F72").Value = "Jonas Ullmann"
$ws.Range("G72").Value = 2

$ws.Range("B73").Value = 2
$ws.Range("F73").Value = "Eva Greven"
# G73 stays at 1 (unchanged)

# Row 74 loses its "Micha Landoll" label entirely, leaving only G74.
$ws.Range("F74").ClearContents()
$ws.Range("G74").Value = 0
